$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ProjectDebt (row 18) for years 1-6 is reset to 0 (it previously ramped
# 1000, 2000, ..., 6000). Years 7-16 already just carry G18 forward via
# formula, so they cascade automatically.
$ws.Range("B18:G18").Value = 0

# Move the active selection to H18, matching the author's final cursor
# position after making the edit.
$ws.Range("H18").Select()
